$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.891.51"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "2.241.30"
$ws.Range("E3").Value = "  +1.92%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "272.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +15.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.648"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0956"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +20.86%  "

$ws.Range("E13").Value = "  +1.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.01%  "

$ws.Range("D15").Value = "2.573.77"
$ws.Range("E15").Value = "  +1.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.823"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.32%  "

$ws.Range("D17").Value = "2.263.42"
$ws.Range("E17").Value = "  +2.87%  "

$ws.Range("D18").Value = "43.881.53"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("E19").Value = "  +2.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.90%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.55%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0939"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.15%  "

$ws.Range("E35").Value = "  +2.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0353"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +23.49%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.235"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +17.91%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.74%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.100"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.61%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("E48").Value = "  +4.75%  "

$ws.Range("E49").Value = "  +2.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.453"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.62%  "

$ws.Range("D51").Value = "2.458.15"
$ws.Range("E51").Value = "  +1.90%  "
